# Update countries & provincias Spain
# Refreshes the COVID country table on sheet "Pais" with a newer data pull
# (timestamp 01:32 -> 02:49) and re-ranks four countries whose totals moved
# them past their neighbour in the sorted-by-cases list. Because the sheet
# is sorted by "Casos totales" descending, promoting a country into its
# neighbour's row pushes the neighbour (and, for multi-place jumps, the
# countries between) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = 'Datos actualizados a 8 de Agosto de 2020 a las 02:49'

# --- Plain data refreshes (no re-ranking) ---------------------------------
# Estados Unidos
$ws.Range("B4").Value = 5095340
$ws.Range("C4").Value = 63062
$ws.Range("D4").Value = 2616891
$ws.Range("E4").Value = 2314355
$ws.Range("G4").Value = 1290
$ws.Range("H4").Value = 164094

# Argentina
$ws.Range("B21").Value = 235677
$ws.Range("C21").Value = 7482
$ws.Range("E21").Value = 127969
$ws.Range("G21").Value = 160
$ws.Range("H21").Value = 4411

# Canada
$ws.Range("B27").Value = 118985
$ws.Range("C27").Value = 424
$ws.Range("D27").Value = 103435
$ws.Range("E27").Value = 6580

# Australia
$ws.Range("D72").Value = 11320
$ws.Range("E72").Value = 8686

# Libia
$ws.Range("B104").Value = 5079
$ws.Range("C104").Value = 200
$ws.Range("D104").Value = 660
$ws.Range("E104").Value = 4311
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 108

# Comoras
$ws.Range("D169").Value = 353
$ws.Range("E169").Value = 36

# Seychelles
$ws.Range("D189").Value = 125
$ws.Range("E189").Value = 1

# Polinesia Francesa
$ws.Range("B194").Value = 69
$ws.Range("C194").Value = 5
$ws.Range("E194").Value = 7

# San Vicente y las Granadinas
$ws.Range("D195").Value = 49
$ws.Range("E195").Value = 7

# --- Paraguay overtakes Mauritania (rows 97-98) ---------------------------
$ws.Range("A97").Value = 'Paraguay'
$ws.Range("B97").Value = 6508
$ws.Range("C97").Value = 133
$ws.Range("D97").Value = 5123
$ws.Range("E97").Value = 1316
$ws.Range("G97").Value = 3
$ws.Range("H97").Value = 69

$ws.Range("A98").Value = 'Mauritania'
$ws.Range("B98").Value = 6498
$ws.Range("C98").Value = 25
$ws.Range("D98").Value = 5443
$ws.Range("E98").Value = 898
$ws.Range("H98").Value = 157

# --- Surinam jumps ahead of Lituania/Estonia/Ruanda (rows 128-131) -------
$ws.Range("A128").Value = 'Surinam'
$ws.Range("B128").Value = 2203
$ws.Range("C128").Value = 107
$ws.Range("D128").Value = 1505
$ws.Range("E128").Value = 669
$ws.Range("H128").Value = 29

$ws.Range("A129").Value = 'Lituania'
$ws.Range("B129").Value = 2194
$ws.Range("C129").Value = 23
$ws.Range("D129").Value = 1658
$ws.Range("E129").Value = 455
$ws.Range("H129").Value = 81

$ws.Range("A130").Value = 'Estonia'
$ws.Range("B130").Value = 2133
$ws.Range("C130").Value = 9
$ws.Range("D130").Value = 1956
$ws.Range("E130").Value = 114
$ws.Range("H130").Value = 63

$ws.Range("A131").Value = 'Ruanda'
$ws.Range("B131").Value = 2128
$ws.Range("C131").Value = 17
$ws.Range("D131").Value = 1297
$ws.Range("E131").Value = 826
$ws.Range("H131").Value = 5

# --- Martinica jumps ahead of Isla de Man/Mongolia/Islas Feroe/Guadalupe/
#     Eritrea (rows 172-177) ----------------------------------------------
$ws.Range("A172").Value = 'Martinica'
$ws.Range("C172").Value = 60
$ws.Range("D172").Value = 98
$ws.Range("E172").Value = 222
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = 16

$ws.Range("A173").Value = 'Isla de Man'
$ws.Range("B173").Value = 336
$ws.Range("D173").Value = 312
$ws.Range("E173").Value = 0
$ws.Range("H173").Value = 24

$ws.Range("A174").Value = 'Mongolia'
$ws.Range("B174").Value = 293
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 260
$ws.Range("E174").Value = 33

$ws.Range("A175").Value = 'Islas Feroe'
$ws.Range("B175").Value = 291
$ws.Range("C175").Value = 25
$ws.Range("D175").Value = 192
$ws.Range("E175").Value = 99
$ws.Range("H175").Value = 0

$ws.Range("A176").Value = 'Guadalupe'
$ws.Range("B176").Value = 290
$ws.Range("C176").Value = 11
$ws.Range("D176").Value = 186
$ws.Range("E176").Value = 90
$ws.Range("H176").Value = 14

$ws.Range("A177").Value = 'Eritrea'
$ws.Range("B177").Value = 285
$ws.Range("C177").Value = 3
$ws.Range("D177").Value = 245
$ws.Range("E177").Value = 40
$ws.Range("H177").Value = 0

# --- Santa Lucia swaps ahead of Timor Oriental (rows 202-203) ------------
# Totals are identical between the two, so only the names swap.
$ws.Range("A202").Value = 'Santa Lucia'
$ws.Range("A203").Value = 'Timor Oriental'
